$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 813.8072129084443
$ws.Range("C2").Value = 1841.465874636106
$ws.Range("D2").Value = 881.7182273524012
